$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.673.41"
$ws.Range("E2").Value = "  +0.87%  "

$ws.Range("D3").Value = "2.237.83"
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("E4").Value = "  +0.47%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.23"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.56"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.566"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.77%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.507"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.41"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.37%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0795"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.07"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.17%  "

$ws.Range("E13").Value = "  +0.23%  "

$ws.Range("D14").Value = "2.583.28"
$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").Value = "2.331.56"
$ws.Range("E15").Value = "  +0.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.822"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.41"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.55%  "

$ws.Range("D18").Value = "44.609.59"
$ws.Range("E18").Value = "  +1.69%  "

$ws.Range("D19").Value = "0.0₃0925"
$ws.Range("E19").Value = "  -3.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.11"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.57"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.13"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.45"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.91"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.04%  "

$ws.Range("E25").Value = "  -2.91%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("E27").Value = "  +3.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.65"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.75"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -6.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.79"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.73"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.83%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "148.49"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.10%  "

$ws.Range("E33").Value = "  +0.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0766"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.11"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.49%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.108"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.43%  "

$ws.Range("E37").Value = "  -2.38%  "

$ws.Range("E38").Value = "  +4.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.79"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.30"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.71"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0294"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.44%  "

$ws.Range("E43").Value = "  +0.68%  "

$ws.Range("D44").Value = "1.796.92"
$ws.Range("E44").Value = "  +3.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.74"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "80.66"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.185"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "96.97"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.31%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.79"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.71%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.12"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.28"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.60%  "
